$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")
$ws2 = $wb.Worksheets.Item("OwnerTotals")

$ws.Range("G4").Value = "Final/OT"
$ws.Range("O4").Value = 44
$ws.Range("G9").Value = "Final/OT"
$ws.Range("O9").Value = 22
$ws.Range("G11").Value = "1:43 - 1st Half"
$ws.Range("H11").Value = 11
$ws.Range("I11").Value = 8
$ws.Range("J11").Value = 3
$ws.Range("O11").Value = 15
$ws.Range("G12").Value = "1:43 - 1st Half"
$ws.Range("H12").Value = 7
$ws.Range("I12").Value = 10
$ws.Range("J12").Value = 1
$ws.Range("O12").Value = 16
$ws.Range("G13").Value = "Final/OT"
$ws.Range("O13").Value = 36
$ws.Range("G17").Value = "1:43 - 1st Half"
$ws.Range("H17").Value = 1
$ws.Range("J17").Value = 1
$ws.Range("K17").Value = 1
$ws.Range("O17").Value = 6
$ws.Range("G18").Value = "Final/OT"
$ws.Range("H18").Value = 30
$ws.Range("I18").Value = 25
$ws.Range("K18").Value = 4
$ws.Range("G19").Value = "1:43 - 1st Half"
$ws.Range("H19").Value = 26
$ws.Range("I19").Value = 20
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 2
$ws.Range("O19").Value = 17
$ws.Range("G21").Value = "Final/OT"
$ws.Range("G25").Value = "Final/OT"
$ws.Range("O25").Value = 41
$ws.Range("G27").Value = "1:43 - 1st Half"
$ws.Range("H27").Value = 3
$ws.Range("J27").Value = 2
$ws.Range("O27").Value = 15
$ws.Range("G28").Value = "1:43 - 1st Half"
$ws.Range("J28").Value = 3
$ws.Range("N28").Value = 3
$ws.Range("O28").Value = 16
$ws.Range("G31").Value = "Final/OT"
$ws.Range("O31").Value = 22
$ws.Range("G32").Value = "Final/OT"
$ws.Range("G34").Value = "Final/OT"
$ws.Range("H34").Value = 15
$ws.Range("I34").Value = 19
$ws.Range("O34").Value = 38
$ws.Range("G38").Value = "Final/OT"
$ws.Range("D39").Value = "Kevin Overton"
$ws.Range("E39").Value = "AUB"
$ws.Range("G39").Value = "1:43 - 1st Half"
$ws.Range("I39").Value = 4
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("N39").Value = 0
$ws.Range("D40").Value = "Xaivian Lee"
$ws.Range("E40").Value = "FLA"
$ws.Range("G40").Value = "1:43 - 1st Half"
$ws.Range("H40").Value = 4
$ws.Range("I40").Value = 3
$ws.Range("J40").Value = 1
$ws.Range("K40").Value = 1
$ws.Range("N40").Value = 2
$ws.Range("O40").Value = 9
$ws.Range("G43").Value = "1:43 - 1st Half"
$ws.Range("G49").Value = "Final/OT"
$ws.Range("G53").Value = "Final/OT"
$ws.Range("G56").Value = "1:43 - 1st Half"
$ws.Range("J56").Value = 3
$ws.Range("O56").Value = 16
$ws.Range("G57").Value = "1:43 - 1st Half"
$ws.Range("H57").Value = 2
$ws.Range("I57").Value = 2
$ws.Range("J57").Value = 1
$ws.Range("O57").Value = 13
$ws.Range("G60").Value = "Final/OT"
$ws.Range("H60").Value = 22
$ws.Range("I60").Value = 21
$ws.Range("G64").Value = "Final/OT"
$ws.Range("D84").Value = "Urban Klavzar"
$ws.Range("E84").Value = "FLA"
$ws.Range("F84").Value = "AUB@FLA"
$ws.Range("G84").Value = "1:43 - 1st Half"
$ws.Range("H84").Value = 5
$ws.Range("I84").Value = 7
$ws.Range("J84").Value = 0
$ws.Range("M84").Value = 0
$ws.Range("N84").Value = 0
$ws.Range("O84").Value = 12
$ws.Range("D85").Value = "Brandon Garrison"
$ws.Range("E85").Value = "UK"
$ws.Range("F85").Value = "MISS@UK"
$ws.Range("J85").Value = 5
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 1
$ws.Range("N85").Value = 1
$ws.Range("O85").Value = 13
$ws.Range("D86").Value = "Simeon Wilcher"
$ws.Range("E86").Value = "TEX"
$ws.Range("F86").Value = "UGA@TEX"
$ws.Range("H86").Value = 4
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 4
$ws.Range("L86").Value = 2
$ws.Range("N86").Value = 0
$ws.Range("O86").Value = 15
$ws.Range("D87").Value = "Dellquan Warren"
$ws.Range("E87").Value = "MSST"
$ws.Range("F87").Value = "VAN@MSST"
$ws.Range("I87").Value = 2
$ws.Range("K87").Value = 2
$ws.Range("N87").Value = 2
$ws.Range("O87").Value = 9
$ws.Range("D88").Value = "Josh Holloway"
$ws.Range("E88").Value = "TA&M"
$ws.Range("F88").Value = "SC@TA&M"
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 4
$ws.Range("L88").Value = 1
$ws.Range("N88").Value = 1
$ws.Range("O88").Value = 12
$ws.Range("D89").Value = "Kareem Stagg"
$ws.Range("E89").Value = "UGA"
$ws.Range("F89").Value = "UGA@TEX"
$ws.Range("I89").Value = 0
$ws.Range("N89").Value = 0
$ws.Range("O89").Value = 13
$ws.Range("D90").Value = "Patton Pinkins"
$ws.Range("E90").Value = "MISS"
$ws.Range("F90").Value = "MISS@UK"
$ws.Range("G90").Value = "Final"
$ws.Range("I90").Value = 9
$ws.Range("J90").Value = 4
$ws.Range("N90").Value = 1
$ws.Range("O90").Value = 30
$ws.Range("D91").Value = "Sebastian Mack"
$ws.Range("E91").Value = "MIZ"
$ws.Range("F91").Value = "OU@MIZ"
$ws.Range("G91").Value = "Final/OT"
$ws.Range("H91").Value = 3
$ws.Range("I91").Value = 4
$ws.Range("J91").Value = 0
$ws.Range("O91").Value = 4
$ws.Range("D92").Value = "Jamie Vinson"
$ws.Range("E92").Value = "TA&M"
$ws.Range("F92").Value = "SC@TA&M"
$ws.Range("I92").Value = 2
$ws.Range("J92").Value = 2
$ws.Range("O92").Value = 5
$ws.Range("D93").Value = "Kanon Catchings"
$ws.Range("E93").Value = "UGA"
$ws.Range("F93").Value = "UGA@TEX"
$ws.Range("G93").Value = "Final"
$ws.Range("O93").Value = 17
$ws.Range("G98").Value = "1:43 - 1st Half"
$ws.Range("O98").Value = 9
$ws.Range("G99").Value = "Final/OT"
$ws.Range("G100").Value = "1:43 - 1st Half"
$ws.Range("G104").Value = "1:43 - 1st Half"
$ws.Range("G108").Value = "Final/OT"
$ws.Range("G110").Value = "Final/OT"
$ws.Range("G112").Value = "1:43 - 1st Half"
$ws.Range("O112").Value = 9
$ws.Range("G114").Value = "Final/OT"
$ws.Range("G115").Value = "Final/OT"
$ws.Range("G116").Value = "Final/OT"

$ws2.Range("A2").Value = "Clay"
$ws2.Range("B2").Value = 69
$ws2.Range("C2").Value = 3
$ws2.Range("A3").Value = "Tar"
$ws2.Range("B3").Value = 63
$ws2.Range("C3").Value = 4
$ws2.Range("B6").Value = 49
$ws2.Range("B7").Value = 45
$ws2.Range("B8").Value = 32
